$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.022680305178376
$ws.Range("D2").Value = 1.026441591250012
$ws.Range("E2").Value = 1.047240088899647
$ws.Range("F2").Value = 1.051236601709801
$ws.Range("I2").Value = 1.028107959857691
$ws.Range("J2").Value = 1.027864780605469
$ws.Range("K2").Value = 1.029264227892524
$ws.Range("L2").Value = 1.050003280226983
$ws.Range("M2").Value = 1.05398866255884
$ws.Range("N2").Value = 1.013297654950742

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.023814574245496
$ws.Range("D3").Value = 1.02723671351503
$ws.Range("E3").Value = 1.048609489330848
$ws.Range("F3").Value = 1.052703932155035
$ws.Range("I3").Value = 1.028268440723941
$ws.Range("J3").Value = 1.028636243807575
$ws.Range("K3").Value = 1.029866898756318
$ws.Range("L3").Value = 1.051183039048338
$ws.Range("M3").Value = 1.055266915320932
$ws.Range("N3").Value = 1.013557109867162

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.024547947313807
$ws.Range("D4").Value = 1.027750536646682
$ws.Range("E4").Value = 1.049495810182668
$ws.Range("F4").Value = 1.053653670127505
$ws.Range("I4").Value = 1.028370636855177
$ws.Range("J4").Value = 1.029134357520505
$ws.Range("K4").Value = 1.030255543978371
$ws.Range("L4").Value = 1.051946117719041
$ws.Range("M4").Value = 1.056093795263164
$ws.Range("N4").Value = 1.013724521103148

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.024856121606375
$ws.Range("D5").Value = 1.027966386836554
$ws.Range("E5").Value = 1.049868477446053
$ws.Range("F5").Value = 1.054053010126463
$ws.Range("I5").Value = 1.028413206181727
$ws.Range("J5").Value = 1.029343508305167
$ws.Range("K5").Value = 1.030418614185663
$ws.Range("L5").Value = 1.052266846349969
$ws.Range("M5").Value = 1.056441362628932
$ws.Range("N5").Value = 1.013794787706435

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.024907857520443
$ws.Range("D6").Value = 1.028002619598504
$ws.Range("E6").Value = 1.049931053423327
$ws.Range("F6").Value = 1.054120065294592
$ws.Range("I6").Value = 1.028420330660633
$ws.Range("J6").Value = 1.029378610647234
$ws.Range("K6").Value = 1.030445975861141
$ws.Range("L6").Value = 1.052320694106832
$ws.Range("M6").Value = 1.056499717709561
$ws.Range("N6").Value = 1.013806579154823

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.02455206568484
$ws.Range("D7").Value = 1.027753421479924
$ws.Range("E7").Value = 1.049500789548597
$ws.Range("F7").Value = 1.053659005847679
$ws.Range("I7").Value = 1.028371207216398
$ws.Range("J7").Value = 1.029137153209027
$ws.Range("K7").Value = 1.030257724172466
$ws.Range("L7").Value = 1.05195040358039
$ws.Range("M7").Value = 1.056098439679292
$ws.Range("N7").Value = 1.013725460452947

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.023063756882457
$ws.Range("D8").Value = 1.026710446333349
$ws.Range("E8").Value = 1.047702839301213
$ws.Range("F8").Value = 1.051732437842109
$ws.Range("I8").Value = 1.028162535913408
$ws.Range("J8").Value = 1.028125723326348
$ws.Range("K8").Value = 1.029468177607172
$ws.Range("L8").Value = 1.050402050139371
$ws.Range("M8").Value = 1.054420704733002
$ws.Range("N8").Value = 1.013385437143701

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.020436679165079
$ws.Range("D9").Value = 1.02486741174451
$ws.Range("E9").Value = 1.044536200834785
$ws.Range("F9").Value = 1.048339529374495
$ws.Range("I9").Value = 1.027782222997511
$ws.Range("J9").Value = 1.02633517786823
$ws.Range("K9").Value = 1.028066738318393
$ws.Range("L9").Value = 1.047671181583713
$ws.Range("M9").Value = 1.051462372557879
$ws.Range("N9").Value = 1.012782634535301

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.018682155814457
$ws.Range("D10").Value = 1.023635214230741
$ws.Range("E10").Value = 1.042425944712062
$ws.Range("F10").Value = 1.046078660631057
$ws.Range("I10").Value = 1.027520198696557
$ws.Range("J10").Value = 1.025135847950305
$ws.Range("K10").Value = 1.027125580340942
$ws.Range("L10").Value = 1.04584874016507
$ws.Range("M10").Value = 1.049488631644772
$ws.Range("N10").Value = 1.012378302154682

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.017921655702976
$ws.Range("D11").Value = 1.023100820384002
$ws.Range("E11").Value = 1.041512328154464
$ws.Range("F11").Value = 1.045099879211584
$ws.Range("I11").Value = 1.027404725624863
$ws.Range("J11").Value = 1.024615172843715
$ws.Range("K11").Value = 1.026716411429079
$ws.Range("L11").Value = 1.045059116166579
$ws.Range("M11").Value = 1.048633572607713
$ws.Range("N11").Value = 1.012202633033729

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.017639052424699
$ws.Range("D12").Value = 1.022902195191224
$ws.Range("E12").Value = 1.041172986627263
$ws.Range("F12").Value = 1.04473633997369
$ws.Range("I12").Value = 1.027361530863232
$ws.Range("J12").Value = 1.024421565525977
$ws.Range("K12").Value = 1.026564180429656
$ws.Range("L12").Value = 1.044765736366598
$ws.Range("M12").Value = 1.048315898954677
$ws.Range("N12").Value = 1.012137292567714

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.017699677243071
$ws.Range("D13").Value = 1.022944806751331
$ws.Range("E13").Value = 1.041245775868625
$ws.Range("F13").Value = 1.044814319375195
$ws.Range("I13").Value = 1.027370809997592
$ws.Range("J13").Value = 1.024463104271435
$ws.Range("K13").Value = 1.026596845707232
$ws.Range("L13").Value = 1.044828670935813
$ws.Range("M13").Value = 1.048384044114308
$ws.Range("N13").Value = 1.012151312359587

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017898298083691
$ws.Range("D14").Value = 1.023084404569118
$ws.Range("E14").Value = 1.041484277770302
$ws.Range("F14").Value = 1.04506982848748
$ws.Range("I14").Value = 1.027401161313168
$ws.Range("J14").Value = 1.024599173401106
$ws.Range("K14").Value = 1.026703833020044
$ws.Range("L14").Value = 1.045034866920725
$ws.Range("M14").Value = 1.048607314988115
$ws.Range("N14").Value = 1.012197233791107

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.01802065900571
$ws.Range("D15").Value = 1.023170398472806
$ws.Range("E15").Value = 1.041631228696004
$ws.Range("F15").Value = 1.045227259091197
$ws.Range("I15").Value = 1.027419821613822
$ws.Range("J15").Value = 1.024682982820922
$ws.Range("K15").Value = 1.026769718619398
$ws.Range("L15").Value = 1.045161900580089
$ws.Range("M15").Value = 1.048744870586217
$ws.Range("N15").Value = 1.012225515677301

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.018732611056933
$ws.Range("D16").Value = 1.023670662345255
$ws.Range("E16").Value = 1.042486580969237
$ws.Range("F16").Value = 1.04614362272268
$ws.Range("I16").Value = 1.027527819795609
$ws.Range("J16").Value = 1.025170374719083
$ws.Range("K16").Value = 1.02715270088055
$ws.Range("L16").Value = 1.04590113414716
$ws.Range("M16").Value = 1.049545369976304
$ws.Range("N16").Value = 1.012389948258782

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.019178989200689
$ws.Range("D17").Value = 1.023984238155806
$ws.Range("E17").Value = 1.043023154747273
$ws.Range("F17").Value = 1.046718480861971
$ws.Range("I17").Value = 1.027595024625099
$ws.Range("J17").Value = 1.025475738344477
$ws.Range("K17").Value = 1.027392495508602
$ws.Range("L17").Value = 1.046364700571141
$ws.Range("M17").Value = 1.050047387759386
$ws.Range("N17").Value = 1.012492934081326

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.019439278937046
$ws.Range("D18").Value = 1.024167060215742
$ws.Range("E18").Value = 1.043336142792937
$ws.Range("F18").Value = 1.047053804507148
$ws.Range("I18").Value = 1.027634029678979
$ws.Range("J18").Value = 1.025653720916956
$ws.Range("K18").Value = 1.027532205355074
$ws.Range("L18").Value = 1.046635043543928
$ws.Range("M18").Value = 1.050340166359311
$ws.Range("N18").Value = 1.012552947019343

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.019528018351206
$ws.Range("D19").Value = 1.024229384026848
$ws.Range("E19").Value = 1.043442866076201
$ws.Range("F19").Value = 1.04716814450637
$ws.Range("I19").Value = 1.027647296429663
$ws.Range("J19").Value = 1.025714386222632
$ws.Range("K19").Value = 1.027579815977394
$ws.Range("L19").Value = 1.046727215577854
$ws.Range("M19").Value = 1.050439989656874
$ws.Range("N19").Value = 1.012573400233114

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.019131104834226
$ws.Range("D20").Value = 1.023950602863898
$ws.Range("E20").Value = 1.042965584079757
$ws.Range("F20").Value = 1.046656802131919
$ws.Range("I20").Value = 1.027587834290676
$ws.Range("J20").Value = 1.025442989263675
$ws.Range("K20").Value = 1.027366784194456
$ws.Range("L20").Value = 1.046314969220075
$ws.Range("M20").Value = 1.049993530156573
$ws.Range("N20").Value = 1.012481890570704

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017839812520401
$ws.Range("D21").Value = 1.023043300017088
$ws.Range("E21").Value = 1.041414044489335
$ws.Range("F21").Value = 1.044994586804467
$ws.Range("I21").Value = 1.027392231971095
$ws.Range("J21").Value = 1.0245591101158
$ws.Range("K21").Value = 1.02667233476073
$ws.Range("L21").Value = 1.044974149535371
$ws.Range("M21").Value = 1.048541569153873
$ws.Range("N21").Value = 1.012183713537918

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.017027232138184
$ws.Range("D22").Value = 1.022472104538482
$ws.Range("E22").Value = 1.040438622307851
$ws.Range("F22").Value = 1.043949621176901
$ws.Range("I22").Value = 1.027267496309254
$ws.Range("J22").Value = 1.024002191094733
$ws.Range("K22").Value = 1.026234274751846
$ws.Range("L22").Value = 1.044130668064642
$ws.Range("M22").Value = 1.047628275944181
$ws.Range("N22").Value = 1.01199572190064

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.017458062962286
$ws.Range("D23").Value = 1.022774976273898
$ws.Range("E23").Value = 1.040955704793389
$ws.Range("F23").Value = 1.044503566004695
$ws.Range("I23").Value = 1.027333787254974
$ws.Range("J23").Value = 1.024297537615525
$ws.Range("K23").Value = 1.026466634686508
$ws.Range("L23").Value = 1.044577857993479
$ws.Range("M23").Value = 1.048112468085966
$ws.Range("N23").Value = 1.012095428864827

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.019152741949598
$ws.Range("D24").Value = 1.023965801455729
$ws.Range("E24").Value = 1.042991597743125
$ws.Range("F24").Value = 1.046684672036456
$ws.Range("I24").Value = 1.027591083893803
$ws.Range("J24").Value = 1.025457787567525
$ws.Range("K24").Value = 1.027378402519405
$ws.Range("L24").Value = 1.046337440822725
$ws.Range("M24").Value = 1.050017866211241
$ws.Range("N24").Value = 1.012486880833622

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.021116385941707
$ws.Range("D25").Value = 1.025344496711721
$ws.Range("E25").Value = 1.045354688810469
$ws.Range("F25").Value = 1.049216471189418
$ws.Range("I25").Value = 1.027882037632481
$ws.Range("J25").Value = 1.0267990646057
$ws.Range("K25").Value = 1.028430251325522
$ws.Range("L25").Value = 1.048377490104028
$ws.Range("M25").Value = 1.052227424786063
$ws.Range("N25").Value = 1.012938906422133

